$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4083.1667
$ws.Range("J40").Value = 4083.1667
$ws.Range("L40").Value = 4083.1667
$ws.Range("N40").Value = -4433.1667
$ws.Range("H53").Value = 71872.86
$ws.Range("I53").Value = 167.875
$ws.Range("J53").Value = 167479.5
$ws.Range("K53").Value = 167.875
$ws.Range("L53").Value = 167479.5
$ws.Range("M53").Value = 469.125
$ws.Range("N53").Value = -168753.5
$ws.Range("H106").Value = 66521.25
$ws.Range("I106").Value = 4320.8335
$ws.Range("J106").Value = 253122.5
$ws.Range("K106").Value = 4320.8335
$ws.Range("L106").Value = 253122.5
$ws.Range("M106").Value = -3689.8335
$ws.Range("N106").Value = -254384.5
$ws.Range("H107").Value = 243
$ws.Range("I107").Value = 273
$ws.Range("J107").Value = 203
$ws.Range("K107").Value = 273
$ws.Range("L107").Value = 203
$ws.Range("M107").Value = 1647
$ws.Range("N107").Value = -4043
$ws.Range("H123").Value = 57699.668
$ws.Range("J123").Value = 57699.668
$ws.Range("L123").Value = 57699.668
$ws.Range("N123").Value = -67499.66800000001
$ws.Range("H126").Value = 75000
$ws.Range("J126").Value = 75000
$ws.Range("L126").Value = 75000
$ws.Range("N126").Value = -84880
$ws.Range("H137").Value = 2334.7708
$ws.Range("I137").Value = 2097.7878
$ws.Range("J137").Value = 2856.1333
$ws.Range("K137").Value = 6293.3634
$ws.Range("L137").Value = 8568.3999
$ws.Range("M137").Value = -3743.3634
$ws.Range("N137").Value = -13668.3999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 4581.6665
$ws.Range("I5").Value = 176.42857
$ws.Range("K5").Value = 176.42857
$ws.Range("M5").Value = -64.42857000000001
$ws.Range("H32").Value = 4572.378
$ws.Range("I32").Value = 4449.023
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 4449.023
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -4162.023
$ws.Range("N32").Value = -10574
$ws.Range("H45").Value = 2208.3333
$ws.Range("I45").Value = 1725
$ws.Range("K45").Value = 1725
$ws.Range("M45").Value = -1348
$ws.Range("H63").Value = 63162444
$ws.Range("I63").Value = 111113110
$ws.Range("J63").Value = 20006840
$ws.Range("K63").Value = 111113110
$ws.Range("L63").Value = 20006840
$ws.Range("M63").Value = -111112424
$ws.Range("N63").Value = -20008212
$ws.Range("H66").Value = 63162444
$ws.Range("I66").Value = 111113110
$ws.Range("J66").Value = 20006840
$ws.Range("K66").Value = 555565550
$ws.Range("L66").Value = 100034200
$ws.Range("M66").Value = -555562118
$ws.Range("N66").Value = -100041064
$ws.Range("H74").Value = 31252842
$ws.Range("I74").Value = 40001556
$ws.Range("K74").Value = 40001556
$ws.Range("M74").Value = -40000682
$ws.Range("H77").Value = 31252842
$ws.Range("I77").Value = 40001556
$ws.Range("K77").Value = 200007780
$ws.Range("M77").Value = -200003412
$ws.Range("H102").Value = 2134843
$ws.Range("I102").Value = 2731870.5
$ws.Range("K102").Value = 2731870.5
$ws.Range("M102").Value = -2730248.5
$ws.Range("H110").Value = 90911460
$ws.Range("I110").Value = 142859090
$ws.Range("J110").Value = 3099.5
$ws.Range("K110").Value = 142859090
$ws.Range("L110").Value = 3099.5
$ws.Range("M110").Value = -142857045
$ws.Range("N110").Value = -7189.5
$ws.Range("H122").Value = 2437.4856
$ws.Range("I122").Value = 1577.2
$ws.Range("J122").Value = 3584.5334
$ws.Range("K122").Value = 4731.6
$ws.Range("L122").Value = 10753.6002
$ws.Range("M122").Value = -2281.6
$ws.Range("N122").Value = -15653.6002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 4581.6665
$ws.Range("I4").Value = 176.42857
$ws.Range("K4").Value = 176.42857
$ws.Range("M4").Value = -61.42857000000001
$ws.Range("H99").Value = 1438
$ws.Range("I99").Value = 1318.3334
$ws.Range("J99").Value = 1489.2858
$ws.Range("K99").Value = 1318.3334
$ws.Range("L99").Value = 1489.2858
$ws.Range("M99").Value = 179.6666
$ws.Range("N99").Value = -4485.2858
$ws.Range("H134").Value = 3529.889
$ws.Range("I134").Value = 2922.5
$ws.Range("K134").Value = 8767.5
$ws.Range("M134").Value = -6232.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 736.8
$ws.Range("I105").Value = 713.0833
$ws.Range("K105").Value = 713.0833
$ws.Range("M105").Value = 1033.9167

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 648.6667
$ws.Range("I2").Value = 545.0526
$ws.Range("K2").Value = 3270.3156
$ws.Range("M2").Value = -3157.3156
$ws.Range("H6").Value = 58.666668
$ws.Range("I6").Value = 43
$ws.Range("K6").Value = 129
$ws.Range("M6").Value = -16
$ws.Range("H59").Value = 30
$ws.Range("I59").Value = 20
$ws.Range("K59").Value = 60
$ws.Range("M59").Value = 480
$ws.Range("H68").Value = 667403.2
$ws.Range("I68").Value = 667403.2
$ws.Range("K68").Value = 2002209.6
$ws.Range("M68").Value = -2001398.6
$ws.Range("H71").Value = 667403.2
$ws.Range("I71").Value = 667403.2
$ws.Range("K71").Value = 6006628.8
$ws.Range("M71").Value = -6002572.8
$ws.Range("H80").Value = 3569.8572
$ws.Range("I80").Value = 4002
$ws.Range("J80").Value = 3397
$ws.Range("K80").Value = 12006
$ws.Range("L80").Value = 10191
$ws.Range("M80").Value = -11070
$ws.Range("N80").Value = -12063
$ws.Range("H83").Value = 3569.8572
$ws.Range("I83").Value = 4002
$ws.Range("J83").Value = 3397
$ws.Range("K83").Value = 36018
$ws.Range("L83").Value = 30573
$ws.Range("M83").Value = -31338
$ws.Range("N83").Value = -39933
$ws.Range("H99").Value = 5888.778
$ws.Range("I99").Value = 4999
$ws.Range("K99").Value = 14997
$ws.Range("M99").Value = -12751
$ws.Range("H122").Value = 40000240
$ws.Range("I122").Value = 300
$ws.Range("J122").Value = 66666864
$ws.Range("K122").Value = 2700
$ws.Range("L122").Value = 600001776
$ws.Range("M122").Value = -250
$ws.Range("N122").Value = -600006676

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1720.25
$ws.Range("I2").Value = 1240.5
$ws.Range("J2").Value = 2200
$ws.Range("K2").Value = 1240.5
$ws.Range("L2").Value = 2200
$ws.Range("M2").Value = -1127.5
$ws.Range("N2").Value = -2426

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2467
$ws.Range("I22").Value = 2467
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2467
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = -2172
$ws.Range("H27").Value = 2467
$ws.Range("I27").Value = 2467
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 2467
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = -2360
$ws.Range("H46").Value = 2685.44
$ws.Range("I46").Value = 2174.875
$ws.Range("J46").Value = 2925.7058
$ws.Range("K46").Value = 2174.875
$ws.Range("L46").Value = 2925.7058
$ws.Range("M46").Value = -1986.875
$ws.Range("N46").Value = -3301.7058
$ws.Range("H61").Value = 2819.353
$ws.Range("J61").Value = 3874.75
$ws.Range("L61").Value = 3874.75
$ws.Range("N61").Value = -4278.75
$ws.Range("H93").Value = 111113120
$ws.Range("I93").Value = 333333340
$ws.Range("J93").Value = 3004
$ws.Range("K93").Value = 333333340
$ws.Range("L93").Value = 3004
$ws.Range("M93").Value = -333332092
$ws.Range("N93").Value = -5500
$ws.Range("H113").Value = 2819.353
$ws.Range("J113").Value = 3874.75
$ws.Range("L113").Value = 3874.75
$ws.Range("N113").Value = -8214.75
$ws.Range("H136").Value = 3872.5818
$ws.Range("I136").Value = 3463.9788
$ws.Range("J136").Value = 6273.125
$ws.Range("K136").Value = 10391.9364
$ws.Range("L136").Value = 18819.375
$ws.Range("M136").Value = -7841.936399999999
$ws.Range("N136").Value = -23919.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2652724.8
$ws.Range("I62").Value = 7941674.5
$ws.Range("J62").Value = 8250
$ws.Range("K62").Value = 7941674.5
$ws.Range("L62").Value = 8250
$ws.Range("M62").Value = -7941050.5
$ws.Range("N62").Value = -9498
$ws.Range("H64").Value = 147777
$ws.Range("I64").Value = 147777
$ws.Range("K64").Value = 147777
$ws.Range("M64").Value = -147529
$ws.Range("H65").Value = 2652724.8
$ws.Range("I65").Value = 7941674.5
$ws.Range("J65").Value = 8250
$ws.Range("K65").Value = 39708372.5
$ws.Range("L65").Value = 41250
$ws.Range("M65").Value = -39705252.5
$ws.Range("N65").Value = -47490
$ws.Range("H67").Value = 147777
$ws.Range("I67").Value = 147777
$ws.Range("K67").Value = 147777
$ws.Range("M67").Value = -146919
$ws.Range("H100").Value = 15629674
$ws.Range("I100").Value = 19236342
$ws.Range("K100").Value = 38472684
$ws.Range("M100").Value = -38472143
$ws.Range("H136").Value = 8046.7393
$ws.Range("I136").Value = 10313.667
$ws.Range("J136").Value = 3796.25
$ws.Range("K136").Value = 30941.001
$ws.Range("L136").Value = 11388.75
$ws.Range("M136").Value = -28391.001
$ws.Range("N136").Value = -16488.75
